$d = $word.ActiveDocument

# Paragraph text is currently "Version 2." split across runs:
#   "Versi" | "on" | " 2" | "." (with proofErr + bookmark interleaved)
# Target text is "Version 1." laid out as:
#   "Version" | " 1." (proofErr + bookmark preserved, trailing "." run removed)

# Work from the end of the paragraph backwards so earlier offsets stay valid.

# 1) Remove the trailing "." run (positions 9-10) -> merges away entirely.
$rPeriod = $d.Range(9, 10)
$rPeriod.Text = ""

# 2) Change " 2" (positions 7-9) to " 1."
$rNum = $d.Range(7, 9)
$rNum.Text = " 1."

# 3) Merge "Versi" + "on" (positions 0-7) into a single "Version" run.
#    A same-text assignment is treated as a no-op by the engine, so first
#    write a differing value to force the merge, then correct it back.
$rWord = $d.Range(0, 7)
$rWord.Text = "Versionn"
$rWord2 = $d.Range(0, 8)
$rWord2.Text = "Version"
